# Remove multiline/trailing-space comments in Frontend_Table_Description.xlsx
# Strips a single trailing space or a single trailing newline from a handful
# of COLUMN_DESCRIPTION cells in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C43").Value = "Formular angelegt von"
$ws.Range("C70").Value = "Formular angelegt von"

$ws.Range("C44").Value = "Formular zuletzt bearbeitet von"
$ws.Range("C71").Value = "Formular zuletzt bearbeitet von"

$ws.Range("C46").Value = "ID Medikationsanalyse"

$ws.Range("C62").Value = "Zeitaufwand Medikationsanalyse"

$ws.Range("C78").Value = "Anderer Hinweisgeber"

$ws.Range("C146").Value = "1. Bewertung von"

$ws.Range("C148").Value = "Zuordnung Meda -> rMRP"

$ws.Range("C184").Value = "2. Bewertung von"

$ws.Range("C226").Value = "ent. Ern."

$ws.Range("C240").Value = "Hb" + [char]0x2193

$wb.Save()
